$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking price strings are not
# auto-converted to numbers by Excel (they must stay inline/shared strings).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.871.28'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '3.331.38'
$ws.Range("E3").Value = '  +1.12%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '519.89'
$ws.Range("E5").Value = '  -0.74%  '
$ws.Range("D6").Value = '172.19'
$ws.Range("E6").Value = '  -5.60%  '
$ws.Range("D7").Value = '0.590'
$ws.Range("E7").Value = '  -2.45%  '
$ws.Range("D8").Value = '3.327.65'
$ws.Range("E8").Value = '  +1.15%  '
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").Value = '0.602'
$ws.Range("E10").Value = '  -3.29%  '
$ws.Range("D11").Value = '52.79'
$ws.Range("E11").Value = '  -11.43%  '
$ws.Range("E12").Value = '  -1.04%  '
$ws.Range("D13").Value = '0.0000255'
$ws.Range("E13").Value = '  -1.50%  '
$ws.Range("D14").Value = '8.95'
$ws.Range("E14").Value = '  -2.37%  '
$ws.Range("D15").Value = '3.879.00'
$ws.Range("E15").Value = '  +2.37%  '
$ws.Range("D16").Value = '3.336.61'
$ws.Range("E16").Value = '  +1.91%  '
$ws.Range("E17").Value = '  -1.41%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '17.42'
$ws.Range("E18").Value = '  -1.60%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '63.592.43'
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").Value = '11.15'
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("D21").Value = '0.953'
$ws.Range("E21").Value = '  -0.86%  '
$ws.Range("D22").Value = '372.76'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '4.22'
$ws.Range("E23").Value = '  +7.00%  '
$ws.Range("E24").Value = '  -0.94%  '
$ws.Range("D25").Value = '81.29'
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("E26").Value = '  -2.95%  '
$ws.Range("D27").Value = '6.17'
$ws.Range("E27").Value = '  +1.45%  '
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").Value = '11.17'
$ws.Range("E29").Value = '  -3.60%  '
$ws.Range("D30").Value = '8.13'
$ws.Range("E30").Value = '  -3.80%  '
$ws.Range("D31").Value = '28.68'
$ws.Range("E31").Value = '  -0.70%  '
$ws.Range("D32").Value = '624.48'
$ws.Range("E32").Value = '  -2.64%  '
$ws.Range("D33").Value = '6.37'
$ws.Range("E33").Value = '  -8.35%  '
$ws.Range("D34").Value = '11.12'
$ws.Range("E34").Value = '  -2.32%  '
$ws.Range("E35").Value = '  -1.63%  '
$ws.Range("D36").Value = '57.69'
$ws.Range("E36").Value = '  -2.53%  '
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").Value = '35.67'
$ws.Range("E38").Value = '  -3.60%  '
$ws.Range("D39").Value = '0.375'
$ws.Range("E39").Value = '  -7.12%  '
$ws.Range("E40").Value = '  +7.84%  '
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("D42").Value = '2.63'
$ws.Range("E42").Value = '  +5.67%  '
$ws.Range("D43").Value = '2.914.52'
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("D44").Value = '0.123'
$ws.Range("E44").Value = '  -3.47%  '
$ws.Range("D45").Value = '3.00'
$ws.Range("E45").Value = '  +3.38%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '0.0393'
$ws.Range("E47").Value = '  -1.77%  '
$ws.Range("E48").Value = '  -5.86%  '
$ws.Range("D49").Value = '2.97'
$ws.Range("E49").Value = '  +0.65%  '
$ws.Range("E50").Value = '  -1.75%  '
$ws.Range("D51").Value = '135.73'
$ws.Range("E51").Value = '  +2.78%  '
